# "0305 Coaches Suspension Updated"
# Remove five coaches whose suspensions were lifted / entered in error:
#   Bilinski, Brian   (row 10)
#   Carroll,  Jason   (row 14)
#   Kennedy,  William (row 29)
#   Lisinski, Tom     (row 30)
#   Marr,     Kristin (row 31)
# and clear the stray "??" note that was next to Torres, Jeremy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so earlier row numbers stay valid as we go.
$ws.Rows(31).Delete() | Out-Null
$ws.Rows(30).Delete() | Out-Null
$ws.Rows(29).Delete() | Out-Null
$ws.Rows(14).Delete() | Out-Null
$ws.Rows(10).Delete() | Out-Null

# After the five deletions, the Torres/Jeremy row (originally row 41) is now
# row 36; clear its leftover "??" comment but keep the cell's formatting.
$ws.Range("E36").Value = ""

# Match the author's final selection/active cell.
$ws.Range("E36").Select() | Out-Null
